# Fix up the "具有相當價值之財產" (property of considerable value) sheet (sheet5):
# the sheet currently has a "broken" header row (row 1 literally duplicates the
# first data row's values instead of containing real column headers), and is
# missing the standard metadata columns (property_category / category / date /
# legislator_name / legislator_id / source_file / index) that the other
# worksheets in this workbook already carry.
#
# This script rewrites row 1 with real headers, adds a "quantity" column, and
# appends the missing metadata columns to both the header row and the existing
# data row - matching the layout already used on the other sheets (e.g. the
# "存款" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

# ---- Row 1: real column headers for the columns that already existed ----
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "quantity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"

# Cell to copy the header formatting (bold + border) from for the brand new
# header cells.
$headerFmt = $ws.Cells.Item(1, 2)

$ws.Cells.Item(1, 6).Value = "property_category"
$headerFmt.Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)

$ws.Cells.Item(1, 7).Value = "category"
$headerFmt.Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)

$ws.Cells.Item(1, 8).Value = "date"
$headerFmt.Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)

$ws.Cells.Item(1, 9).Value = "legislator_name"
$headerFmt.Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)

$ws.Cells.Item(1, 10).Value = "legislator_id"
$headerFmt.Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)

$ws.Cells.Item(1, 11).Value = "source_file"
$headerFmt.Copy()
$ws.Cells.Item(1, 11).PasteSpecial(-4122)

$ws.Cells.Item(1, 12).Value = "index"
$headerFmt.Copy()
$ws.Cells.Item(1, 12).PasteSpecial(-4122)

# ---- Row 2: existing data row gains the metadata columns ----
# Cell to copy the plain data formatting from for the brand new data cells.
$dataFmt = $ws.Cells.Item(2, 2)

$ws.Cells.Item(2, 6).Value = "otherbonds"
$dataFmt.Copy()
$ws.Cells.Item(2, 6).PasteSpecial(-4122)

$ws.Cells.Item(2, 7).Value = "normal"
$dataFmt.Copy()
$ws.Cells.Item(2, 7).PasteSpecial(-4122)

# Store the date as plain text (matching the other sheets) instead of
# letting Excel auto-convert the "2012-04-30" literal into a date serial.
$dateCell = $ws.Cells.Item(2, 8)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2012-04-30"
$dataFmt.Copy()
$dateCell.PasteSpecial(-4122)

$ws.Cells.Item(2, 9).Value = "蘇震清"
$dataFmt.Copy()
$ws.Cells.Item(2, 9).PasteSpecial(-4122)

$ws.Cells.Item(2, 10).Value = 1718
$dataFmt.Copy()
$ws.Cells.Item(2, 10).PasteSpecial(-4122)

$ws.Cells.Item(2, 11).Value = "tmp16a71"
$dataFmt.Copy()
$ws.Cells.Item(2, 11).PasteSpecial(-4122)

$ws.Cells.Item(2, 12).Value = 79
$dataFmt.Copy()
$ws.Cells.Item(2, 12).PasteSpecial(-4122)
